$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 428.22223
$ws.Range("I2").Value = 261
$ws.Range("K2").Value = 261
$ws.Range("M2").Value = -148
$ws.Range("H16").Value = 57010
$ws.Range("J16").Value = 57010
$ws.Range("L16").Value = 57010
$ws.Range("N16").Value = -57470
$ws.Range("H21").Value = 36279.332
$ws.Range("J21").Value = 9819
$ws.Range("L21").Value = 9819
$ws.Range("N21").Value = -10755
$ws.Range("H23").Value = 36279.332
$ws.Range("J23").Value = 9819
$ws.Range("L23").Value = 9819
$ws.Range("N23").Value = -10287
$ws.Range("H38").Value = 1405
$ws.Range("I38").Value = 1405
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4215
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -3843
$ws.Range("H51").Value = 2687.5
$ws.Range("J51").Value = 4750
$ws.Range("L51").Value = 4750
$ws.Range("N51").Value = -5718
$ws.Range("H58").Value = 2017
$ws.Range("J58").Value = 2017
$ws.Range("L58").Value = 6051
$ws.Range("N58").Value = -6351
$ws.Range("H88").Value = 2951.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2951.5
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 2951.5
$ws.Range("N88").Value = -3763.5
$ws.Range("H91").Value = 2951.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2951.5
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 2951.5
$ws.Range("N91").Value = -5759.5
$ws.Range("H107").Value = 911.8333
$ws.Range("I107").Value = 474
$ws.Range("J107").Value = 1787.5
$ws.Range("K107").Value = 474
$ws.Range("L107").Value = 1787.5
$ws.Range("M107").Value = 1446
$ws.Range("N107").Value = -5627.5
$ws.Range("H116").Value = 10314.917
$ws.Range("J116").Value = 2189.25
$ws.Range("L116").Value = 2189.25
$ws.Range("N116").Value = -9073.25
$ws.Range("H125").Value = 555.44446
$ws.Range("I125").Value = 555.44446
$ws.Range("K125").Value = 4999.00014
$ws.Range("M125").Value = -2539.00014
$ws.Range("H138").Value = 3126.9792
$ws.Range("I138").Value = 4622.75
$ws.Range("J138").Value = 2379.0938
$ws.Range("K138").Value = 13868.25
$ws.Range("L138").Value = 7137.2814
$ws.Range("M138").Value = -8728.25
$ws.Range("N138").Value = -17417.2814
$ws.Range("H139").Value = 71799.5
$ws.Range("J139").Value = 71799.5
$ws.Range("L139").Value = 71799.5
$ws.Range("N139").Value = -82079.5
$ws.Range("H140").Value = 81516.336
$ws.Range("J140").Value = 81516.336
$ws.Range("L140").Value = 81516.336
$ws.Range("N140").Value = -91876.336

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 347670.8
$ws.Range("I2").Value = 463385.1
$ws.Range("K2").Value = 463385.1
$ws.Range("M2").Value = -463272.1
$ws.Range("H32").Value = 4134.922
$ws.Range("I32").Value = 3658.049
$ws.Range("K32").Value = 3658.049
$ws.Range("M32").Value = -3371.049
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0
$ws.Range("H45").Value = 1688.3
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 1835.375
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 1835.375
$ws.Range("M45").Value = -723
$ws.Range("N45").Value = -2589.375
$ws.Range("H61").Value = 6750.35
$ws.Range("I61").Value = 7073.8667
$ws.Range("K61").Value = 7073.8667
$ws.Range("M61").Value = -6861.8667
$ws.Range("H74").Value = 1292.1578
$ws.Range("I74").Value = 606.3333
$ws.Range("K74").Value = 606.3333
$ws.Range("M74").Value = 267.6667
$ws.Range("H77").Value = 1292.1578
$ws.Range("I77").Value = 606.3333
$ws.Range("K77").Value = 3031.6665
$ws.Range("M77").Value = 1336.3335
$ws.Range("H116").Value = 347670.8
$ws.Range("I116").Value = 463385.1
$ws.Range("K116").Value = 463385.1
$ws.Range("M116").Value = -461091.1
$ws.Range("H136").Value = 6750.35
$ws.Range("I136").Value = 7073.8667
$ws.Range("K136").Value = 21221.6001
$ws.Range("M136").Value = -18671.6001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 347670.8
$ws.Range("I3").Value = 463385.1
$ws.Range("K3").Value = 463385.1
$ws.Range("M3").Value = -463271.1
$ws.Range("H86").Value = 103241.15
$ws.Range("I86").Value = 3430.7646
$ws.Range("K86").Value = 3430.7646
$ws.Range("M86").Value = -2307.7646
$ws.Range("H89").Value = 103241.15
$ws.Range("I89").Value = 3430.7646
$ws.Range("K89").Value = 17153.823
$ws.Range("M89").Value = -11537.823
$ws.Range("H134").Value = 7824.4116
$ws.Range("I134").Value = 8791.286
$ws.Range("K134").Value = 26373.858
$ws.Range("M134").Value = -23838.858
$ws.Range("H138").Value = 61997.145
$ws.Range("J138").Value = 61997.145
$ws.Range("L138").Value = 61997.145
$ws.Range("N138").Value = -72277.14499999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H134").Value = 2449.3914
$ws.Range("I134").Value = 2206.9
$ws.Range("J134").Value = 4066
$ws.Range("K134").Value = 6620.700000000001
$ws.Range("L134").Value = 12198
$ws.Range("M134").Value = -4085.700000000001
$ws.Range("N134").Value = -17268

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 960.3333
$ws.Range("I98").Value = 439
$ws.Range("J98").Value = 1064.6
$ws.Range("K98").Value = 1317
$ws.Range("L98").Value = 3193.8
$ws.Range("M98").Value = 181
$ws.Range("N98").Value = -6189.799999999999
$ws.Range("H131").Value = 24884
$ws.Range("J131").Value = 26670.592
$ws.Range("L131").Value = 80011.776
$ws.Range("N131").Value = -90091.776
$ws.Range("H137").Value = 3661.9167
$ws.Range("I137").Value = 1770.2667
$ws.Range("J137").Value = 6814.6665
$ws.Range("K137").Value = 5310.800099999999
$ws.Range("L137").Value = 20443.9995
$ws.Range("M137").Value = -210.8000999999995
$ws.Range("N137").Value = -30643.9995

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2918.6667
$ws.Range("I102").Value = 2629.1
$ws.Range("K102").Value = 2629.1
$ws.Range("M102").Value = -1007.1
$ws.Range("H132").Value = 1540486.9
$ws.Range("I132").Value = 2025879.5
$ws.Range("J132").Value = 3410.3333
$ws.Range("K132").Value = 6077638.5
$ws.Range("L132").Value = 10230.9999
$ws.Range("M132").Value = -6075108.5
$ws.Range("N132").Value = -15290.9999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 810.5
$ws.Range("I55").Value = 826.6667
$ws.Range("J55").Value = 800.8
$ws.Range("K55").Value = 826.6667
$ws.Range("L55").Value = 800.8
$ws.Range("M55").Value = -653.6667
$ws.Range("N55").Value = -1146.8
$ws.Range("H68").Value = 1578
$ws.Range("I68").Value = 1578
$ws.Range("K68").Value = 1578
$ws.Range("M68").Value = -829
$ws.Range("H71").Value = 1578
$ws.Range("I71").Value = 1578
$ws.Range("K71").Value = 7890
$ws.Range("M71").Value = -4146
$ws.Range("H82").Value = 2144.2
$ws.Range("I82").Value = 1900
$ws.Range("J82").Value = 2248.8572
$ws.Range("K82").Value = 1900
$ws.Range("L82").Value = 2248.8572
$ws.Range("M82").Value = -1539
$ws.Range("N82").Value = -2970.8572
$ws.Range("H85").Value = 2144.2
$ws.Range("I85").Value = 1900
$ws.Range("J85").Value = 2248.8572
$ws.Range("K85").Value = 1900
$ws.Range("L85").Value = 2248.8572
$ws.Range("M85").Value = -652
$ws.Range("N85").Value = -4744.8572
$ws.Range("H108").Value = 25626
$ws.Range("J108").Value = 25626
$ws.Range("L108").Value = 25626
$ws.Range("N108").Value = -33306
$ws.Range("H136").Value = 2268.4375
$ws.Range("I136").Value = 1945.909
$ws.Range("K136").Value = 5837.727000000001
$ws.Range("M136").Value = -3287.727000000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 37109
$ws.Range("J70").Value = 37109
$ws.Range("L70").Value = 37109
$ws.Range("N70").Value = -37739
$ws.Range("H73").Value = 37109
$ws.Range("J73").Value = 37109
$ws.Range("L73").Value = 37109
$ws.Range("N73").Value = -39293
$ws.Range("H136").Value = 27779462
$ws.Range("I136").Value = 50506280
$ws.Range("K136").Value = 151518840
$ws.Range("M136").Value = -151516290
